# Auto-applies the cryptos-list refresh described by the commit diff.
# D-column (Price) cells are forced to remain TEXT (matching the original
# inlineStr storage) by briefly flipping NumberFormat to "@" before the
# assignment, then resetting the style back to "Normal" so no stray
# number-format/style is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "68.368.19"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.05%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.645.75"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "597.39"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "158.79"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.82%  "
$ws.Range("E8").Value = "  -0.92%  "
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("E15").Value = "  -2.27%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "68.360.93"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.20%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.626.86"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.28%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "11.38"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.07%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "359.79"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.07%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.41"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.18%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.41"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("E22").Value = "  -2.35%  "
$ws.Range("E23").Value = "  +0.45%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "74.36"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.78"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  -2.49%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.55%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "561.19"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("E33").Value = "  +1.14%  "
$ws.Range("E34").Value = "  +3.83%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -1.42%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "160.25"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "19.66"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.53%  "
$ws.Range("E39").Value = "  -0.98%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.86"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("E42").Value = "  -1.40%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.0₆0323"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -4.13%  "
$ws.Range("E44").Value = "  +0.03%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "157.48"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("E46").Value = "  +0.87%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "21.99"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("E49").Value = "  -1.94%  "
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("E51").Value = "  -0.18%  "
